$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# --- Rows 2-22: price (D) and volume (E) updates ---
Set-TextValue 2 4 "29.353.49"
Set-TextValue 2 5 "  +0.62%  "
Set-TextValue 3 4 "1.936.92"
Set-TextValue 3 5 "  +1.72%  "
Set-TextValue 4 4 "1.003"
Set-TextValue 4 5 "  -0.02%  "
Set-TextValue 5 4 "325.41"
Set-TextValue 5 5 "  -0.16%  "
Set-TextValue 6 4 "1.001"
Set-TextValue 6 5 "  -0.03%  "
Set-TextValue 7 4 "0.4620"
Set-TextValue 7 5 "  +0.28%  "
Set-TextValue 8 4 "0.3872"
Set-TextValue 8 5 "  -0.46%  "
Set-TextValue 9 4 "45.89"
Set-TextValue 9 5 "  -0.36%  "
Set-TextValue 10 4 "0.07825"
Set-TextValue 10 5 "  -0.74%  "
Set-TextValue 11 4 "0.9743"
Set-TextValue 11 5 "  -1.72%  "
Set-TextValue 12 4 "22.60"
Set-TextValue 12 5 "  +3.09%  "
Set-TextValue 13 4 "1.930.95"
Set-TextValue 13 5 "  +3.45%  "
Set-TextValue 14 4 "7.081"
Set-TextValue 14 5 "  +0.35%  "
Set-TextValue 15 4 "5.760"
Set-TextValue 15 5 "  -0.36%  "
Set-TextValue 16 4 "0.07025"
Set-TextValue 16 5 "  -0.07%  "
Set-TextValue 17 4 "86.60"
Set-TextValue 17 5 "  -1.59%  "
Set-TextValue 18 4 "1.004"
Set-TextValue 18 5 "  -0.03%  "
Set-TextValue 19 4 "0.000009808"
Set-TextValue 19 5 "  -1.29%  "
Set-TextValue 20 4 "17.08"
Set-TextValue 20 5 "  +0.08%  "
Set-TextValue 21 4 "1.001"
Set-TextValue 21 5 "  -0.02%  "
Set-TextValue 22 4 "29.398.27"
Set-TextValue 22 5 "  +0.75%  "

# --- Rows 23-51: new BitDAO row inserted, all following rows shift down one; 
#     last row (Quant) drops off the bottom of the A1:E51 range ---
$ws.Cells.Item(23, 2).Value = "BitDAO"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
Set-TextValue 23 4 "0.5010"
Set-TextValue 23 5 "  +1.17%  "
$ws.Cells.Item(24, 2).Value = "Uniswap"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue 24 4 "5.483"
Set-TextValue 24 5 "  +3.08%  "
$ws.Cells.Item(25, 2).Value = "Cosmos"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue 25 4 "11.04"
Set-TextValue 25 5 "  -0.95%  "
$ws.Cells.Item(26, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue 26 4 "2.173.96"
Set-TextValue 26 5 "  +3.56%  "
$ws.Cells.Item(27, 2).Value = "Toncoin"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue 27 4 "2.095"
Set-TextValue 27 5 "  -0.42%  "
$ws.Cells.Item(28, 2).Value = "Monero"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue 28 4 "157.32"
Set-TextValue 28 5 "  +0.76%  "
$ws.Cells.Item(29, 2).Value = "EthereumClassic"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue 29 4 "19.38"
Set-TextValue 29 5 "  -0.54%  "
$ws.Cells.Item(30, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue 30 4 "5.741"
Set-TextValue 30 5 "  -2.90%  "
$ws.Cells.Item(31, 2).Value = "BitcoinCash"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue 31 4 "118.38"
Set-TextValue 31 5 "  -0.28%  "
$ws.Cells.Item(32, 2).Value = "LidoDAOToken"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue 32 4 "1.860"
Set-TextValue 32 5 "  -0.65%  "
$ws.Cells.Item(33, 2).Value = "Stellar"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue 33 4 "0.09354"
Set-TextValue 33 5 "  +0.05%  "
$ws.Cells.Item(34, 2).Value = "ImmutableX"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue 34 4 "0.8591"
Set-TextValue 34 5 "  -3.89%  "
$ws.Cells.Item(35, 2).Value = "Filecoin"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue 35 4 "5.169"
Set-TextValue 35 5 "  -1.56%  "
$ws.Cells.Item(36, 2).Value = "ARBITRUM"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue 36 4 "1.306"
Set-TextValue 36 5 "  -1.03%  "
$ws.Cells.Item(37, 2).Value = "HuobiToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue 37 4 "3.094"
Set-TextValue 37 5 "  -2.05%  "
$ws.Cells.Item(38, 2).Value = "Hedera"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue 38 4 "0.05763"
Set-TextValue 38 5 "  -0.50%  "
$ws.Cells.Item(39, 2).Value = "TrustWalletToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue 39 4 "1.156"
Set-TextValue 39 5 "  -1.35%  "
$ws.Cells.Item(40, 2).Value = "VeChain"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue 40 4 "0.02078"
Set-TextValue 40 5 "  -0.59%  "
$ws.Cells.Item(41, 2).Value = "FraxShare"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue 41 4 "7.657"
Set-TextValue 41 5 "  -0.08%  "
$ws.Cells.Item(42, 2).Value = "TheSandbox"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue 42 4 "0.5671"
Set-TextValue 42 5 "  -0.38%  "
$ws.Cells.Item(43, 2).Value = "Algorand"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue 43 4 "0.1776"
Set-TextValue 43 5 "  -1.20%  "
$ws.Cells.Item(44, 2).Value = "Aptos"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue 44 4 "9.416"
Set-TextValue 44 5 "  -3.18%  "
$ws.Cells.Item(45, 2).Value = "PEPE"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue 45 4 "0.000002910"
Set-TextValue 45 5 "  +40.97%  "
$ws.Cells.Item(46, 2).Value = "MXToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue 46 4 "2.717"
Set-TextValue 46 5 "  +6.32%  "
$ws.Cells.Item(47, 2).Value = "Decentraland"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue 47 4 "0.5288"
Set-TextValue 47 5 "  -1.23%  "
$ws.Cells.Item(48, 2).Value = "EnergySwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue 48 4 "11.41"
Set-TextValue 48 5 "  -4.30%  "
$ws.Cells.Item(49, 2).Value = "Cronos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue 49 4 "0.06866"
Set-TextValue 49 5 "  -2.14%  "
$ws.Cells.Item(50, 2).Value = "RenderToken"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 50 4 "2.077"
Set-TextValue 50 5 "  -6.07%  "
$ws.Cells.Item(51, 2).Value = "NEARProtocol"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue 51 4 "1.817"
Set-TextValue 51 5 "  -1.68%  "
